$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.232.76"
$ws.Range("E2").Value = "  +2.33%  "

$ws.Range("D3").Value = "2.533.06"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.48%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("D9").Value = "2.532.60"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("E11").Value = "  +2.63%  "

$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("E15").Value = "  +2.25%  "

$ws.Range("D16").Value = "2.993.40"
$ws.Range("E16").Value = "  +0.43%  "

$ws.Range("D17").Value = "68.215.41"
$ws.Range("E17").Value = "  +2.60%  "

$ws.Range("D18").Value = "2.534.58"
$ws.Range("E18").Value = "  +0.60%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "368.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.70%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.28%  "

$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.62%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.77%  "

$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("D28").Value = "2.660.65"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  +1.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "541.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.79%  "

$ws.Range("E32").Value = "  +1.91%  "

$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.02%  "

$ws.Range("E35").Value = "  -0.94%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("E38").Value = "  -0.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.66%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("E42").Value = "  +0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.09%  "

$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "147.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.561"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.01%  "

$ws.Range("E49").Value = "  +2.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.08%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0757"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "

